$p = $ppt.ActivePresentation

# --- Slide 3: simplify "Autoencoders (Baseline Model)" -> "Autoencoders" ---
$s3 = $p.Slides.Item(3)
$shape3 = $s3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange
$tr3.Paragraphs(5).Runs(1).Text = "Autoencoders"

# --- Slide 9: update the captions/labels text box ---
$s9 = $p.Slides.Item(9)
$shape9 = $s9.Shapes.Item(2)

# Resize/reposition the textbox
$shape9.Left = 674.5631
$shape9.Width = 229.2816

$tr9 = $shape9.TextFrame.TextRange

# Paragraph 2 (blank line under "Added noise") -> "PSNR: 12.66"
$para2 = $tr9.Paragraphs(2)
$para2.Text = "PSNR: 12.66"
$para2.Font.Size = 19
$para2.Font.Bold = $true
$para2.Font.Color.RGB = 0x292929
$para2.Font.Name = "source-serif-pro"

# Paragraph 5 (blank line under "Our Denoising") -> "PSNR: 24.31"
$para5 = $tr9.Paragraphs(5)
$para5.Text = "PSNR: 24.31"
$para5.Font.Size = 19
$para5.Font.Bold = $true
$para5.Font.Color.RGB = 0xC07000
$para5.Font.Name = "source-serif-pro"

Write-Host "Done step 1"
